# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# existing header style used by the other header cells (e.g. G1 "sum").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last header cell onto the new header cell, then
# set its text — this reuses the existing bordered/bold/centered style
# instead of fabricating a near-duplicate style entry.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Save values for data rows 2-6.
$saveValues = @(0, 1, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
